$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''68.611.99'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.99%  '
$ws.Range('D3').Value = '''3.866.65'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.27%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('D4').Style = "Normal"
$ws.Range('D5').Value = '''602.95'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +0.58%  '
$ws.Range('D6').Value = '''172.61'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  +3.41%  '
$ws.Range('D7').Value = '''3.865.32'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').Value = '''0.531'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +1.01%  '
$ws.Range('E10').Value = '  +2.55%  '
$ws.Range('D11').Value = '''6.55'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  +3.72%  '
$ws.Range('E12').Value = '  +1.28%  '
$ws.Range('D13').Value = '''0.0000290'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +16.56%  '
$ws.Range('D14').Value = '''37.29'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.87%  '
$ws.Range('D15').Value = '''4.516.88'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.31%  '
$ws.Range('D16').Value = '''3.857.89'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -0.07%  '
$ws.Range('D17').Value = '''68.656.91'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +0.91%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D18').Value = '''7.55'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +1.85%  '
$ws.Range('B19').Value = 'Chainlink'
$ws.Range('C19').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D19').Value = '''18.37'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  +0.72%  '
$ws.Range('E20').Value = '  +0.90%  '
$ws.Range('E21').Value = '  +1.85%  '
$ws.Range('D22').Value = '''472.78'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  +1.03%  '
$ws.Range('D23').Value = '''0.735'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.92%  '
$ws.Range('E24').Value = '  +1.28%  '
$ws.Range('D25').Value = '''84.03'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = '  +0.85%  '
$ws.Range('D26').Value = '''2.29'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +2.88%  '
$ws.Range('D27').Value = '''12.32'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +1.40%  '
$ws.Range('D28').Value = '''10.51'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +5.38%  '
$ws.Range('E29').Value = '  +0.05%  '
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('D31').Value = '''4.019.55'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('D32').Value = '''7.81'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +1.42%  '
$ws.Range('D33').Value = '''2.33'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +1.04%  '
$ws.Range('D34').Value = '''31.34'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  +1.02%  '
$ws.Range('D35').Value = '''9.49'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +0.93%  '
$ws.Range('D36').Value = '''3.833.09'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +0.01%  '
$ws.Range('D37').Value = '''3.98'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +21.25%  '
$ws.Range('E38').Value = '  +1.90%  '
$ws.Range('E39').Value = '  +2.32%  '
$ws.Range('E40').Value = '  +0.59%  '
$ws.Range('E41').Value = '  +0.83%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('E43').Value = '  +2.93%  '
$ws.Range('D44').Value = '''0.000305'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +12.09%  '
$ws.Range('D45').Value = '''2.01'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = '  +1.45%  '
$ws.Range('B46').Value = 'Cosmos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D46').Value = '''8.82'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +3.12%  '
$ws.Range('B47').Value = 'USDe'
$ws.Range('C47').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D47').Value = '''1.00'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.00%  '
$ws.Range('D48').Value = '''422.24'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -1.66%  '
$ws.Range('D49').Value = '''46.79'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -0.91%  '
$ws.Range('D50').Value = '''0.0362'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +2.39%  '
$ws.Range('D51').Value = '''142.33'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  -0.74%  '
